# Adjustments for additional level or requirements levels
#
# The "Zulagen" column header (in sheet "Vorlage_Datenblatt", row 1) currently
# reads:
#   Zulagen
#   (Schicht-, Sonntags- und Nachtarbeit sowie andere Erschwerniszulagen, 1/12 der Jahressumme)
#
# It needs to drop the trailing ", 1/12 der Jahressumme" clause so it reads:
#   Zulagen
#   (Schicht-, Sonntags- und Nachtarbeit sowie andere Erschwerniszulagen)
#
# The cell is rich text: the "Zulagen" run inherits the cell's (bold) style,
# while the parenthetical note is a second run in a smaller, non-bold Arial.
# We trim the text via the Characters() collection (instead of replacing
# .Value wholesale) so the two-run rich-text formatting survives the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vorlage_Datenblatt")
$cell = $ws.Range("L1")

$removed = ", 1/12 der Jahressumme"
$fullText = $cell.Text
$pos = $fullText.IndexOf($removed)

if ($pos -ge 0) {
    # Delete the unwanted fragment, keeping everything else (and its
    # per-run formatting) untouched.
    $cell.Characters($pos + 1, $removed.Length).Text = ""

    # Re-assert the (unchanged) formatting of the note run so the engine
    # keeps it as its own run distinct from the bold "Zulagen" heading.
    $newText = $cell.Text
    $noteStart = "Zulagen".Length + 1
    $noteLen = $newText.Length - "Zulagen".Length
    $note = $cell.Characters($noteStart, $noteLen)
    $note.Font.Name = "Arial"
    $note.Font.Size = 9
    $note.Font.Bold = $false
}

# The row's wrapped/rotated header text reflows slightly once the note
# shrinks; match the resulting row height.
$ws.Rows.Item(1).RowHeight = 218.25
